# DOMA-3097 review fixes for contact export
#
# The "Unit Type" column (with sample value "квартира") used to be the LAST
# column (F) in the contact-import example sheet. Review feedback moved it
# right after "Unit Name", i.e. it is now column C - "Phones", "Full Name"
# and "Email" each shift one column to the right (C->D, D->E, E->F).
#
# Implemented as a column move: cut column F and insert it before column C.
# The mailto hyperlink that lived on the old Email cell (E2) has to be
# re-created on its new location (F2), since moving/inserting columns does
# not repoint existing Hyperlink objects automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Unit Type" column (F) to become the new column C; Phones/Full
# Name/Email shift right by one column (C->D, D->E, E->F).
$ws.Columns("F").Cut()
$ws.Columns("C").Insert()

# The Email cell's hyperlink stayed bound to E2 - recreate it on F2, where
# the Email column (and its "test@example.com" value) now lives.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:test@example.com", "", "", "test@example.com")

# Adding the hyperlink re-styles the cell with the built-in "Hyperlink"
# look (underline + theme color). Restore the plain formatting the Email
# column used everywhere else in the sheet.
$f2 = $ws.Range("F2")
$f2.Font.Name = "Calibri"
$f2.Font.Size = 12
$f2.Font.Underline = $false
$f2.Font.Color = 0
